$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# AVL section (rows 3-7)
$ws.Range("C3").Value = 0.016
$ws.Range("D3").Value = 0.027
$ws.Range("C4").Value = 0.10100000000000001
$ws.Range("D4").Value = 0.186
$ws.Range("C5").Value = 2.9609999999999999
$ws.Range("D5").Value = 4.4880000000000004
$ws.Range("C6").Value = 46.213999999999999
$ws.Range("D6").Value = 63.720999999999997
$ws.Range("C7").Value = 971.01499999999999
$ws.Range("D7").Value = 1076.9690000000001

# std::map section (rows 11-15)
$ws.Range("C11").Value = 0.049000000000000002
$ws.Range("D11").Value = 0.245
$ws.Range("C12").Value = 0.45900000000000002
$ws.Range("D12").Value = 0.503
$ws.Range("C13").Value = 4.4260000000000002
$ws.Range("D13").Value = 4.0730000000000004
$ws.Range("C14").Value = 60.356999999999999
$ws.Range("D14").Value = 45.289000000000001
$ws.Range("C15").Value = 1239.923
$ws.Range("D15").Value = 990.65300000000002

# Log section (rows 19-23)
$ws.Range("C19").Value = 0.042000000000000003
$ws.Range("D19").Value = 0.055
$ws.Range("C20").Value = 0.254
$ws.Range("D20").Value = 0.38600000000000001
$ws.Range("C21").Value = 3.282
$ws.Range("D21").Value = 3.6539999999999999
$ws.Range("C22").Value = 68.432000000000002
$ws.Range("D22").Value = 59.259
$ws.Range("C23").Value = 1465.337
$ws.Range("D23").Value = 1270.432

# Update selection to match the saved cursor position
$ws.Range("C25").Select()
